$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new value for column F (dSF)
$changes = @{
    2  = 0
    7  = -3
    11 = -4
    13 = -3
    14 = 0
    17 = 1
    21 = 3
    31 = 7
    36 = -1
    37 = -3
    38 = -1
    39 = -1
    41 = -2
    53 = -1
    56 = -7
    59 = -4
    61 = 8
    65 = 1
    67 = 0
    72 = -6
    75 = 6
    76 = 2
    77 = -1
    79 = 2
    83 = 2
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
